# Natmi following Dr Hou advice
#
# The Cntn1/Ptprz1 ligand-receptor sheet originally only reported the
# sCs -> sCs cluster pair. Re-run across the FAPs / sCs / ECs cluster
# combinations (3 senders x 3 targets restricted to the 6 pairs that
# actually co-occur) and extend the data rows below the header
# accordingly (existing row 2 is overwritten, rows 3-7 are new).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then the 20 column values (A..T) in order:
#   Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#   Ligand-expressing cells, Ligand detection rate,
#   Ligand average expression value, Ligand total expression value,
#   Ligand derived specificity of average/total expression value,
#   Receptor-expressing cells, Receptor detection rate,
#   Receptor average expression value, Receptor total expression value,
#   Receptor derived specificity of average/total expression value,
#   Edge average/total expression weight,
#   Edge average/total expression derived specificity
$data = @(
    @(2, @("FAPs","Cntn1","Ptprz1","ECs",
           3,1,0.1145166666666667,0.34355,0.5923587361059505,0.6855055141740312,
           1,0.5,0.0191995,0.038399,0.002780311579676609,0.002767389187860387,
           0.002198662741666667,0.01319197645,0.001646941853317975,0.001897060548143889)),
    @(3, @("FAPs","Cntn1","Ptprz1","FAPs",
           3,1,0.1145166666666667,0.34355,0.5923587361059505,0.6855055141740312,
           2,0.6666666666666666,0.06449100000000001,0.193473,0.009339049146327988,0.01394346436998132,
           0.007385294350000001,0.06646764915,0.005532067348750203,0.009558321712311328)),
    @(4, @("FAPs","Cntn1","Ptprz1","sCs",
           3,1,0.1145166666666667,0.34355,0.5923587361059505,0.6855055141740312,
           2,1,6.821830500000001,13.643661,0.9878806392739954,0.9832891464421584,
           0.7812132894250001,4.687279736550001,0.5851797269038823,0.674050131913576)),
    @(5, @("sCs","Cntn1","Ptprz1","ECs",
           1,0.5,0.0788065,0.157613,0.4076412638940496,0.3144944858259688,
           1,0.5,0.0191995,0.038399,0.002780311579676609,0.002767389187860387,
           0.00151304539675,0.006052181587000001,0.001133369726358635,0.0008703286397164978)),
    @(6, @("sCs","Cntn1","Ptprz1","FAPs",
           1,0.5,0.0788065,0.157613,0.4076412638940496,0.3144944858259688,
           2,0.6666666666666666,0.06449100000000001,0.193473,0.009339049146327988,0.01394346436998132,
           0.0050823099915,0.030493859949,0.003806981797577786,0.00438514265766999)),
    @(7, @("sCs","Cntn1","Ptprz1","sCs",
           1,0.5,0.0788065,0.157613,0.4076412638940496,0.3144944858259688,
           2,1,6.821830500000001,13.643661,0.9878806392739954,0.9832891464421584,
           0.5376045852982501,2.150418341193,0.4027009123701132,0.3092390145285823))
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($entry in $data) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
